$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Eg" label to "Whole Eg"
$ws.Range("A2").Value = "Whole Eg"

# Remove the "Pet food" / "Roasted coffee" rows (rows 11-12)
$ws.Rows("11:12").Delete()

# Remove the trailing rows no longer present in the list
# (Propane/kerosene/firewood, Flour and prepared flour mixes,
#  Bacon and related products, Fresh whole chicken) which are now
# rows 13-16 after the deletion above.
$ws.Rows("13:16").Delete()
